$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status for f7c573b5 file (row 3) in both language columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: f7c573b5 row (row 3) is now handed back ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-11-10 06:42:39"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 12.9166666666667

# --- de-de sheet: f7c573b5 row (row 3) is now handed back ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-11-10 06:42:58"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 12.9166666666667
